# Update countries & provincias Spain
#
# The source "paises.xlsx" table (sheet "Pais") lists one country per row
# (columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes).
#
# This refreshes the COVID-19 figures to a newer snapshot. Most countries
# keep their row position, but a handful of countries that were (almost)
# tied on "Casos totales" swap ranking order with their neighboring row,
# so both the country name (column A) and the statistics (columns B-H)
# for that pair of rows are exchanged/updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Iran (row 11) - updated figures, no reordering
$ws.Range('B11').Value = 77995
$ws.Range('C11').Value = 1606
$ws.Range('D11').Value = 52229
$ws.Range('E11').Value = 20897
$ws.Range('F11').Value = 3594
$ws.Range('G11').Value = 92
$ws.Range('H11').Value = 4869

# Barein (row 63) - updated figures, no reordering
$ws.Range('B63').Value = 1698
$ws.Range('C63').Value = 27
$ws.Range('D63').Value = 703
$ws.Range('E63').Value = 988

# Hungria (row 64) - updated figures, no reordering
$ws.Range('F64').Value = 60

# Azerbaiyan / Eslovenia (rows 72-73) swap order, with updated figures
$ws.Range('A72').Value = 'Eslovenia'
$ws.Range('B72').Value = 1268
$ws.Range('C72').Value = 20
$ws.Range('D72').Value = 174
$ws.Range('E72').Value = 1033
$ws.Range('F72').Value = 31
$ws.Range('H72').Value = 61

$ws.Range('A73').Value = 'Azerbaiyan'
$ws.Range('B73').Value = 1253
$ws.Range('D73').Value = 404
$ws.Range('E73').Value = 836
$ws.Range('F73').Value = 24
$ws.Range('H73').Value = 13

# Albania (row 96) - updated figures, no reordering
$ws.Range('B96').Value = 518
$ws.Range('C96').Value = 24
$ws.Range('D96').Value = 277
$ws.Range('E96').Value = 215
$ws.Range('G96').Value = 1
$ws.Range('H96').Value = 26

# Vietnam (row 115) - updated figures, no reordering
$ws.Range('D115').Value = 176
$ws.Range('E115').Value = 92

# Madagascar (row 135) - updated figures, no reordering
$ws.Range('B135').Value = 111
$ws.Range('C135').Value = 1
$ws.Range('D135').Value = 33
$ws.Range('E135').Value = 78

# Bermudas / Togo (rows 142-143) swap order, with updated figures
$ws.Range('A142').Value = 'Togo'
$ws.Range('D142').Value = 35
$ws.Range('F142').Value = 0
$ws.Range('H142').Value = 3

$ws.Range('A143').Value = 'Bermudas'
$ws.Range('D143').Value = 33
$ws.Range('F143').Value = 3
$ws.Range('H143').Value = 5

# San Martin (Parte Holandesa) / Bahamas (rows 154-155) swap order, with updated figures
$ws.Range('A154').Value = 'Bahamas'
$ws.Range('D154').Value = 6
$ws.Range('F154').Value = 1
$ws.Range('H154').Value = 8

$ws.Range('A155').Value = 'San Martin (Parte Holandesa)'
$ws.Range('D155').Value = 5
$ws.Range('F155').Value = 2
$ws.Range('H155').Value = 9

# Nueva Caledonia / Timor Oriental (rows 177-178) swap order, with updated figures
$ws.Range('A177').Value = 'Timor Oriental'
$ws.Range('C177').Value = 10
$ws.Range('F177').Value = 0

$ws.Range('A178').Value = 'Nueva Caledonia'
$ws.Range('C178').Value = 0
$ws.Range('F178').Value = 1

# Islas Virgenes de los Estados Unidos / Fiyi (rows 180-181) swap order, with updated figures
$ws.Range('A180').Value = 'Fiyi'
$ws.Range('C180').Value = 1

$ws.Range('A181').Value = 'Islas Virgenes de los Estados Unidos'
$ws.Range('C181').Value = 0

# Suazilandia / Dominica (rows 185-186) swap order, with updated figures
$ws.Range('A185').Value = 'Dominica'
$ws.Range('C185').Value = 0

$ws.Range('A186').Value = 'Suazilandia'
$ws.Range('C186').Value = 1

# Granada / San Cristobal y Nieves (rows 189-190) swap order, with updated figures
$ws.Range('A189').Value = 'San Cristobal y Nieves'
$ws.Range('F189').Value = 0

$ws.Range('A190').Value = 'Granada'
$ws.Range('F190').Value = 2

# Yemen / San Pedro y Miquelon (rows 215-216) swap order (figures identical, only names change)
$ws.Range('A215').Value = 'San Pedro y Miquelon'
$ws.Range('A216').Value = 'Yemen'
